$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header): extend with two new columns, P1 = 14, Q1 = 15 ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
# Copy formatting from O1 (existing header style) onto the new header cells
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# --- Data rows 2-25: flip I/K/M/O values and append new P/Q columns ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P: new column, value 2
    $ws.Cells.Item($r, 17).Value = 2   # Q: new column, value 2
}
